$wb = $excel.ActiveWorkbook

# "groups" sheet (sheet1): selection moves from G8 to E2
$wsGroups = $wb.Worksheets.Item(1)
$wsGroups.Range("E2").Select()

# "projects" sheet (sheet2): update the name_with_namespace / full_path
# values for the NULL-1 / NULL-2 / NULL-3 rows, then move the selection
# from D10 to F8. This sheet stays the active tab.
$wsProjects = $wb.Worksheets.Item(2)
$wsProjects.Activate()

$wsProjects.Range("E2").Value = "ansible-roles/NULL-1"
$wsProjects.Range("D2").Value = "ansible-roles / NULL-1"
$wsProjects.Range("D3").Value = "ansible-roles / NULL-2"
$wsProjects.Range("E3").Value = "ansible-roles/NULL-2"
$wsProjects.Range("D4").Value = "ansible-roles / NULL-3"
$wsProjects.Range("E4").Value = "ansible-roles/NULL-2"

$wsProjects.Range("F8").Select()
